# Controle dos experimentos NP2 - "Plot 15 MB 8C SP 1C"
# Swap the two annotation blocks on sheet "1 Cliente":
#  - E17 label ("Observar a geração do bitrate que tá estranha") becomes "gerar de novo"
#  - E26 label ("rodou agora") becomes "rodando"
#  - The ID columns below each label (A18:A22 and A27:A31) swap their numbering
#  - B31 (stray helper value) is cleared
#  - The active selection moves from E26 to E20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1 Cliente")

# Update the two text labels. Order matters: write E26 first so the
# now-unreferenced shared-string slots are reclaimed in the same order
# as in the target file (slot 21 -> "rodando", slot 22 -> "gerar de novo").
$ws.Range("E26").Value = "rodando"
$ws.Range("E17").Value = "gerar de novo"

# First block (rows 18-22) now carries the numbering that used to live
# in the second block (7,8,9,10,11)
$ws.Range("A18").Value = 7
$ws.Range("A19").Value = 8
$ws.Range("A20").Value = 9
$ws.Range("A21").Value = 10
$ws.Range("A22").Value = 11

# Second block (rows 27-31) now carries the numbering that used to live
# in the first block (1,2,3,4,5)
$ws.Range("A27").Value = 1
$ws.Range("A28").Value = 2
$ws.Range("A29").Value = 3
$ws.Range("A30").Value = 4
$ws.Range("A31").Value = 5

# Row 31's helper column B is no longer populated
$ws.Range("B31").ClearContents()

# Move the active selection
$ws.Activate()
$ws.Range("E20").Select()
